# haver pull with social security
# gftfbdx was missing in haver pivoted
#
# Adds a new "code"/"reference" pair row to Sheet1:
#   A74 = gftfbdx, B74 = social_security
# and moves the active selection down to the new row, matching the
# author's on-disk view state after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new haver code + friendly name as the next row after the
# existing data (previously the sheet ended at row 73 -> A1:B73).
$ws.Range("A74").Value = "gftfbdx"
$ws.Range("B74").Value = "social_security"

# Keep the sheet's selection state consistent with the new data extent
# (selection follows the row right after the newly appended one, same
# as the author's saved view).
$ws.Range("A75").Select()
